# chore: update Sheets via scheduled runner
# Refreshes market-price-derived leve profit figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) across the crafting-class sheets. Source data
# has no formulas -- every cell below is a plain numeric literal, so the
# "scheduled runner" simply re-pushes the latest fetched values.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 53
$ws.Range("H53").Value = 132.77777
$ws.Range("I53").Value = 69.28570999999999
$ws.Range("J53").Value = 173.18182
$ws.Range("K53").Value = 69.28570999999999
$ws.Range("L53").Value = 173.18182
$ws.Range("M53").Value = 567.71429
$ws.Range("N53").Value = -1447.18182

# Row 86
$ws.Range("H86").Value = 7146377
$ws.Range("I86").Value = 8335519.5
$ws.Range("K86").Value = 8335519.5
$ws.Range("M86").Value = -8334396.5

# Row 89
$ws.Range("H89").Value = 7146377
$ws.Range("I89").Value = 8335519.5
$ws.Range("K89").Value = 41677597.5
$ws.Range("M89").Value = -41671981.5

# Row 94
$ws.Range("H94").Value = 601744.4
$ws.Range("I94").Value = 601744.4
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 601744.4
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -601293.4

# Row 113
$ws.Range("H113").Value = 2235.3845
$ws.Range("I113").Value = 2197.5
$ws.Range("J113").Value = 2239.7144
$ws.Range("K113").Value = 2197.5
$ws.Range("L113").Value = 2239.7144
$ws.Range("M113").Value = 1056.5
$ws.Range("N113").Value = -8747.714400000001

# Row 125
$ws.Range("H125").Value = 1554.6666
$ws.Range("I125").Value = 2032
$ws.Range("J125").Value = 600
$ws.Range("K125").Value = 18288
$ws.Range("L125").Value = 5400
$ws.Range("M125").Value = -15828
$ws.Range("N125").Value = -10320

# Row 137
$ws.Range("H137").Value = 4091.0667
$ws.Range("I137").Value = 1270.2273
$ws.Range("K137").Value = 3810.6819
$ws.Range("M137").Value = -1260.6819

# Row 141
$ws.Range("H141").Value = 1415
$ws.Range("I141").Value = 734.2105
$ws.Range("J141").Value = 6958.5713
$ws.Range("K141").Value = 2202.6315
$ws.Range("L141").Value = 20875.7139
$ws.Range("M141").Value = 2977.3685
$ws.Range("N141").Value = -31235.7139

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 8839.402
$ws.Range("I32").Value = 8313.046
$ws.Range("K32").Value = 8313.046
$ws.Range("M32").Value = -8026.046

# Row 61
$ws.Range("H61").Value = 1219.4423
$ws.Range("I61").Value = 985.14636
$ws.Range("K61").Value = 985.14636
$ws.Range("M61").Value = -773.14636

# Row 122
$ws.Range("H122").Value = 2286.5
$ws.Range("I122").Value = 2372.353
$ws.Range("J122").Value = 1800
$ws.Range("K122").Value = 7117.059
$ws.Range("L122").Value = 5400
$ws.Range("M122").Value = -4667.059
$ws.Range("N122").Value = -10300

# Row 132
$ws.Range("H132").Value = 10001565
$ws.Range("I132").Value = 14286709
$ws.Range("K132").Value = 42860127
$ws.Range("M132").Value = -42857597

# Row 136
$ws.Range("H136").Value = 1219.4423
$ws.Range("I136").Value = 985.14636
$ws.Range("K136").Value = 2955.43908
$ws.Range("M136").Value = -405.4390800000001

# Row 139
$ws.Range("H139").Value = 48671.285
$ws.Range("J139").Value = 48671.285
$ws.Range("L139").Value = 48671.285
$ws.Range("N139").Value = -58951.285

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 112
$ws.Range("H112").Value = 39059
$ws.Range("J112").Value = 39059
$ws.Range("L112").Value = 39059
$ws.Range("N112").Value = -42013

# Row 132
$ws.Range("H132").Value = 35926
$ws.Range("J132").Value = 35926
$ws.Range("L132").Value = 35926
$ws.Range("N132").Value = -46046

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 62
$ws.Range("H62").Value = 3387.3333
$ws.Range("I62").Value = 2950
$ws.Range("J62").Value = 3737.2
$ws.Range("K62").Value = 2950
$ws.Range("L62").Value = 3737.2
$ws.Range("M62").Value = -2326
$ws.Range("N62").Value = -4985.2

# Row 65
$ws.Range("H65").Value = 3387.3333
$ws.Range("I65").Value = 2950
$ws.Range("J65").Value = 3737.2
$ws.Range("K65").Value = 14750
$ws.Range("L65").Value = 18686
$ws.Range("M65").Value = -11630
$ws.Range("N65").Value = -24926

# Row 94
$ws.Range("H94").Value = 2688.889
$ws.Range("I94").Value = 1966.6666
$ws.Range("J94").Value = 3050
$ws.Range("K94").Value = 1966.6666
$ws.Range("L94").Value = 3050
$ws.Range("M94").Value = -1515.6666
$ws.Range("N94").Value = -3952

# Row 99
$ws.Range("H99").Value = 2564.2222
$ws.Range("J99").Value = 2342.3333
$ws.Range("L99").Value = 2342.3333
$ws.Range("N99").Value = -5338.3333

# Row 126
$ws.Range("H126").Value = 2564.2222
$ws.Range("J126").Value = 2342.3333
$ws.Range("L126").Value = 7026.999899999999
$ws.Range("N126").Value = -11966.9999

# Row 134
$ws.Range("H134").Value = 342599.25
$ws.Range("I134").Value = 953.6286
$ws.Range("J134").Value = 2335532
$ws.Range("K134").Value = 2860.8858
$ws.Range("L134").Value = 7006596
$ws.Range("M134").Value = -325.8858
$ws.Range("N134").Value = -7011666

# Row 137
$ws.Range("H137").Value = 66912.86
$ws.Range("J137").Value = 66912.86
$ws.Range("L137").Value = 66912.86
$ws.Range("N137").Value = -77112.86

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 21
$ws.Range("H21").Value = 7800
$ws.Range("J21").Value = 7800
$ws.Range("L21").Value = 7800
$ws.Range("N21").Value = -8146

# Row 30
$ws.Range("H30").Value = 7800
$ws.Range("J30").Value = 7800
$ws.Range("L30").Value = 7800
$ws.Range("N30").Value = -8010

# Row 82
$ws.Range("H82").Value = 28000
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()

# Row 85
$ws.Range("H85").Value = 28000
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()

# Row 122
$ws.Range("H122").Value = 1550
$ws.Range("I122").Value = 1616.6666
$ws.Range("K122").Value = 4849.9998
$ws.Range("M122").Value = -2399.9998

# Row 126
$ws.Range("H126").Value = 23782.4
$ws.Range("I126").Value = 37337.332
$ws.Range("J126").Value = 3450
$ws.Range("K126").Value = 112011.996
$ws.Range("L126").Value = 10350
$ws.Range("M126").Value = -109541.996
$ws.Range("N126").Value = -15290

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 22
$ws.Range("H22").Value = 555.9231
$ws.Range("I22").Value = 485.2857
$ws.Range("K22").Value = 485.2857
$ws.Range("M22").Value = -190.2857

# Row 27
$ws.Range("H27").Value = 555.9231
$ws.Range("I27").Value = 485.2857
$ws.Range("K27").Value = 485.2857
$ws.Range("M27").Value = -378.2857

# Row 132
$ws.Range("H132").Value = 2504.3833
$ws.Range("I132").Value = 1625.8889
$ws.Range("J132").Value = 5139.8667
$ws.Range("K132").Value = 4877.6667
$ws.Range("L132").Value = 15419.6001
$ws.Range("M132").Value = -2347.6667
$ws.Range("N132").Value = -20479.6001

# Row 136
$ws.Range("H136").Value = 1170.1086
$ws.Range("I136").Value = 930.8605
$ws.Range("K136").Value = 2792.5815
$ws.Range("M136").Value = -242.5815000000002

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 128
$ws.Range("H128").Value = 49501.668
$ws.Range("J128").Value = 49501.668
$ws.Range("L128").Value = 49501.668
$ws.Range("N128").Value = -59461.668

# Row 132
$ws.Range("H132").Value = 1914.3914
$ws.Range("I132").Value = 1653.2941
$ws.Range("J132").Value = 2654.1667
$ws.Range("K132").Value = 4959.8823
$ws.Range("L132").Value = 7962.500100000001
$ws.Range("M132").Value = -2429.8823
$ws.Range("N132").Value = -13022.5001

# Row 136
$ws.Range("H136").Value = 181221.52
$ws.Range("I136").Value = 250613.25
$ws.Range("K136").Value = 751839.75
$ws.Range("M136").Value = -749289.75

# Row 139
$ws.Range("H139").Value = 56850
$ws.Range("J139").Value = 56850
$ws.Range("L139").Value = 56850
$ws.Range("N139").Value = -67130
